$d = $word.ActiveDocument

# XML namespace declaration snippet reused for every InsertXML call.
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function New-ParaXml($text, $extra) {
    $escaped = $text -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'
    # Word only emits xml:space="preserve" when the text has leading/trailing
    # whitespace that would otherwise be collapsed.
    if ($text -ne $text.Trim()) {
        $tOpen = 'xml:space="preserve"'
    } else {
        $tOpen = ''
    }
    return "<w:p $wns><w:r><w:t $tOpen>$escaped</w:t></w:r>$extra</w:p>"
}

# 1. "Task 01" -> "Task 02", and move the _GoBack bookmark here (it currently
#    sits at the end of the "RunnableTask task2 = new RunnableTask();" line).
$d.Paragraphs(2).Range.InsertXML((New-ParaXml "Task 02" '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'))

# 2. Re-write the Java source lines as clean single-run paragraphs, which
#    drops all the stray w:proofErr spell/grammar-check markers left over
#    from splitting runs.
$d.Paragraphs(3).Range.InsertXML((New-ParaXml "package multithreadapp;" ''))
$d.Paragraphs(4).Range.InsertXML((New-ParaXml "public class RunnableTask implements Runnable {" ''))
$d.Paragraphs(5).Range.InsertXML((New-ParaXml "    @Override" ''))
$d.Paragraphs(6).Range.InsertXML((New-ParaXml "    public void run() {" ''))
$d.Paragraphs(7).Range.InsertXML((New-ParaXml "        System.out.println(Thread.currentThread().getId() + `" is executing the runnable task.`");" ''))
$d.Paragraphs(8).Range.InsertXML((New-ParaXml "    }" ''))
$d.Paragraphs(9).Range.InsertXML((New-ParaXml "    public static void main(String[] args) {" ''))
$d.Paragraphs(10).Range.InsertXML((New-ParaXml "        RunnableTask task1 = new RunnableTask();" ''))
$d.Paragraphs(11).Range.InsertXML((New-ParaXml "        RunnableTask task2 = new RunnableTask();" ''))
$d.Paragraphs(12).Range.InsertXML((New-ParaXml "        Thread thread1 = new Thread(task1);" ''))
$d.Paragraphs(13).Range.InsertXML((New-ParaXml "        Thread thread2 = new Thread(task2);" ''))
$d.Paragraphs(14).Range.InsertXML((New-ParaXml "        thread1.start(); // Starts thread1" ''))
$d.Paragraphs(15).Range.InsertXML((New-ParaXml "        thread2.start(); // Starts thread2" ''))
